$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-9 from
# 45243 (2023-11-13) to 45244 (2023-11-14), keeping existing formatting.
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45244
}
